# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for the rows
# whose dialog-act classification changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 7; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 16; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 25; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 27; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 28; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 32; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 35; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 37; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 41; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 42; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 52; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 55; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 61; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 69; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 84; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 90; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 96; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 106; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 115; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 130; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 135; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 140; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 171; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 173; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 201; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 213; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 215; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 221; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 226; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 233; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 239; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 240; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 248; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 261; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 270; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 273; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 276; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 293; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 295; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 300; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 308; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 324; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 326; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 346; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 356; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 373; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 380; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}
